# Apply targeted cell value corrections to Sheet1 (X_train_sel_k4 data).
# These are spot-fixes to specific BMI (C), Glucose (B), Pregnancies (A)
# and Age (D) values that were cleared/left as placeholders in the
# original export and now get their finished/imputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 13.35
$ws.Range("C21").Value = 50.55
$ws.Range("C44").Value = 13.35
$ws.Range("C64").Value = 50.55
$ws.Range("B99").Value = 37.125
$ws.Range("D135").Value = 66.5
$ws.Range("C140").Value = 50.55
$ws.Range("B154").Value = 37.125
$ws.Range("D155").Value = 66.5
$ws.Range("C156").Value = 13.35
$ws.Range("C157").Value = 13.35
$ws.Range("C171").Value = 13.35
$ws.Range("B177").Value = 37.125
$ws.Range("A185").Value = 13.5
$ws.Range("C188").Value = 50.55
$ws.Range("B218").Value = 37.125
$ws.Range("C219").Value = 13.35
$ws.Range("C230").Value = 50.55
$ws.Range("D244").Value = 66.5
$ws.Range("A275").Value = 13.5
$ws.Range("D313").Value = 66.5
$ws.Range("C391").Value = 50.55
$ws.Range("C400").Value = 13.35
$ws.Range("D400").Value = 66.5
$ws.Range("D412").Value = 66.5
$ws.Range("A473").Value = 13.5
$ws.Range("B501").Value = 37.125
$ws.Range("D509").Value = 66.5
$ws.Range("A536").Value = 13.5
$ws.Range("D588").Value = 66.5
